$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$values = @(
    "1+79=80",
    "90-24=66",
    "14+46=60",
    "61-49=12",
    "9-6=3",
    "17+52=69",
    "90-89=1",
    "60-36=24",
    "82-3=79",
    "78-58=20",
    "90+9=99",
    "34+18=52",
    "97-72=25",
    "47-38=9",
    "70-10=60",
    "96-46=50",
    "24+41=65",
    "73-44=29",
    "88-23=65",
    "52+37=89",
    "3+71=74",
    "96-5=91",
    "45+41=86",
    "1+55=56",
    "42-20=22",
    "57-18=39",
    "87-11=76",
    "79+9=88",
    "44+10=54",
    "46-34=12",
    "64-61=3",
    "93-33=60",
    "42-14=28",
    "55+31=86",
    "83-82=1",
    "82-65=17",
    "76+10=86",
    "93-0=93",
    "85-39=46",
    "77-16=61",
    "16+58=74",
    "52+3=55",
    "54-1=53",
    "81-31=50",
    "81-19=62",
    "51+44=95",
    "55+44=99",
    "36+6=42",
    "79-7=72",
    "29-26=3",
    "41-2=39",
    "55+35=90",
    "51-11=40",
    "1+16=17",
    "69-10=59",
    "97-85=12",
    "42+5=47",
    "37+21=58",
    "18+80=98",
    "62-22=40",
    "36+14=50",
    "41+47=88",
    "84+0=84",
    "17+34=51",
    "12+87=99",
    "74-30=44",
    "66-29=37",
    "39+15=54",
    "79+13=92",
    "86-18=68",
    "33-5=28",
    "24-2=22",
    "74+23=97",
    "44+38=82",
    "40-10=30",
    "92-5=87",
    "19+31=50",
    "34-7=27",
    "64-36=28",
    "62+34=96",
    "31+68=99",
    "50-37=13",
    "96-25=71",
    "80+9=89",
    "97-25=72",
    "86-68=18",
    "44-21=23",
    "43-25=18",
    "65+20=85",
    "37+19=56",
    "44-2=42",
    "18-0=18",
    "78-57=21",
    "96-26=70",
    "51+17=68",
    "22+62=84",
    "48-45=3",
    "62-54=8",
    "64-9=55",
    "36-5=31"
)

$idx = 0
for ($r = 1; $r -le 20; $r++) {
    for ($c = 1; $c -le 5; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}
Write-Host "Updated" $idx "cells"